{"js": "// Apply the five resume text revisions described by the diff.\n// Each entry is an exact, literal, case-sensitive search string paired\n// with its replacement. We locate the run containing the old text via\n// Body.search() (matchCase, no wildcards so the text is matched\n// literally, spaces and all) and then rewrite it in place with\n// Range.insertText(..., \"Replace\") so all other run formatting\n// (rPr/font/size) is preserved.\nconst replacements = [\n  {\n    old: \" -                                           (STEM) \",\n    new: \" -                                           Systems  \"\n  },\n  {\n    old:\n      \"The working utility invention Natural Human WaveLength & Impedance Meter/Visual Recognition Medical Instrument I engineered and built, submitted a white paper to the U.S Army Research Laboratory(ARL).\",\n    new:\n      \"The working utility invention; Natural Human WaveLength & Impedance Meter/Visual Recognition Medical Instrument I formulated, designed, engineered and built, submitted a white paper to the U.S Army Research Laboratory(ARL).\"\n  },\n  {\n    old:\n      'Engineered and built Toke Core the working hyperprotovisor systems software running the Bitcoin p2p network, confirming and reporting to the U.S. Federal Trade Commission at antitrust@ftc.gov the illegitimacy of cryptocurrency and inefficiency of \"blockchain\" as a hashing log verifier, \"hyperledger\".',\n    new:\n      'Developed, engineered and built; Toke Core, the working hyperprotovisor systems software running the Bitcoin p2p network on github.com/TokeBit, confirming and communicating to the U.S. Federal Trade Commission at antitrust@ftc.gov the illegitimacy of cryptocurrency and inefficiency of \"blockchain\" as a hashing log verifier, \"hyperledger\".'\n  },\n  {\n    old:\n      \"First and only to achieve Nuclear Fusion by implicitly engineering and building the working utility invention Nuclear Fusion Reactor and a working Fusion Reactionary Engine,  reported to the U.S. National Aeronautical & Space Agency(NASA) for geo-satellite concatenation of scientific proofs. \",\n    new:\n      \"First and only to achieve Nuclear Fusion by implicitly formulating, designing, engineering and building the working utility invention; Nuclear Fusion Reactor and Fusion Reactionary Engine, communicated to the U.S. National Aeronautical & Space Agency(NASA) for geo-satellite concatenation of scientific proofs. \"\n  },\n  {\n    old:\n      \"Engineered and partially built the working utility invention Full-Spatial Median-Free Liquid and Photonic Bit Transfer Module System 276,480-bit Computational Processor and submitted a drafted technical manual to the U.S. National Security Agency(NSA).\",\n    new:\n      \"Formulated, designing,  engineered and partially built the working utility invention; Full-Spatial Median-Free Liquid and Photonic Bit Transfer Module System 276,480-bit Computational Processor and submitted a drafted technical manual to the U.S. National Security Agency(NSA).\"\n  }\n];\n\nfor (const { old, new: replacement } of replacements) {\n  const results = context.document.body.search(old, {\n    matchCase: true,\n    matchWholeWord: false\n  });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + old);\n  }\n\n  for (const hit of results.items) {\n    hit.insertText(replacement, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Apply the five resume text revisions described by the diff.\n#\n# For each entry we locate the exact, literal, case-sensitive old text\n# with Find.Execute (MatchWildcards:$false so parentheses/ampersands/etc.\n# are treated literally, not as wildcard metacharacters) and then assign\n# the replacement directly to Range.Text. Using Range.Text (rather than\n# passing the replacement through Find.Execute's ReplaceWith/Replace\n# arguments) keeps the run's existing formatting (rPr) intact and avoids\n# Word's Find-and-Replace AutoCorrect pass that would otherwise turn the\n# straight double quotes in the Toke Core bullet into curly/smart quotes.\n\n$d = $word.ActiveDocument\n\nfunction Replace-ExactText($doc, [string]$oldText, [string]$newText) {\n    $rng = $doc.Content\n    $rng.Find.ClearFormatting()\n    $found = $rng.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, \"\", 0)\n    if (-not $found) {\n        throw \"Text not found: $oldText\"\n    }\n    $rng.Text = $newText\n}\n\nReplace-ExactText $d \" -                                           (STEM) \" \" -                                           Systems  \"\n\nReplace-ExactText $d \"The working utility invention Natural Human WaveLength & Impedance Meter/Visual Recognition Medical Instrument I engineered and built, submitted a white paper to the U.S Army Research Laboratory(ARL).\" \"The working utility invention; Natural Human WaveLength & Impedance Meter/Visual Recognition Medical Instrument I formulated, designed, engineered and built, submitted a white paper to the U.S Army Research Laboratory(ARL).\"\n\nReplace-ExactText $d 'Engineered and built Toke Core the working hyperprotovisor systems software running the Bitcoin p2p network, confirming and reporting to the U.S. Federal Trade Commission at antitrust@ftc.gov the illegitimacy of cryptocurrency and inefficiency of \"blockchain\" as a hashing log verifier, \"hyperledger\".' 'Developed, engineered and built; Toke Core, the working hyperprotovisor systems software running the Bitcoin p2p network on github.com/TokeBit, confirming and communicating to the U.S. Federal Trade Commission at antitrust@ftc.gov the illegitimacy of cryptocurrency and inefficiency of \"blockchain\" as a hashing log verifier, \"hyperledger\".'\n\nReplace-ExactText $d \"First and only to achieve Nuclear Fusion by implicitly engineering and building the working utility invention Nuclear Fusion Reactor and a working Fusion Reactionary Engine,  reported to the U.S. National Aeronautical & Space Agency(NASA) for geo-satellite concatenation of scientific proofs. \" \"First and only to achieve Nuclear Fusion by implicitly formulating, designing, engineering and building the working utility invention; Nuclear Fusion Reactor and Fusion Reactionary Engine, communicated to the U.S. National Aeronautical & Space Agency(NASA) for geo-satellite concatenation of scientific proofs. \"\n\nReplace-ExactText $d \"Engineered and partially built the working utility invention Full-Spatial Median-Free Liquid and Photonic Bit Transfer Module System 276,480-bit Computational Processor and submitted a drafted technical manual to the U.S. National Security Agency(NSA).\" \"Formulated, designing,  engineered and partially built the working utility invention; Full-Spatial Median-Free Liquid and Photonic Bit Transfer Module System 276,480-bit Computational Processor and submitted a drafted technical manual to the U.S. National Security Agency(NSA).\"\n"}
